$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("S2").Value = 0.01611173663836548
$ws.Range("T2").Value = 0.01611173663836548

# Row 3 updates
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 0.2800659258742222
$ws.Range("R3").Value = 2.520593332868001
$ws.Range("S3").Value = 0.68328279700753
$ws.Range("T3").Value = 0.68328279700753

# Row 4 updates
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 0.1232130365728889
$ws.Range("R4").Value = 1.108917329156
$ws.Range("S4").Value = 0.3006054663541045
$ws.Range("T4").Value = 0.3006054663541044
